# Add basic File implementation + open-ended FileMetadata sheet
# (mirrors: project/excel/acr_harmonized_data_model.xlsx diff)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "File" worksheet - appended after the last existing sheet (FamilyMember)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$fileSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$fileSheet.Name = "File"

$fileHeaders = @(
    "subject",
    "sample",
    "filename",
    "format",
    "data_type",
    "size",
    "drs_uri",
    "file_metadata",
    "has_access_policy",
    "id",
    "external_id"
)

for ($i = 0; $i -lt $fileHeaders.Length; $i++) {
    $fileSheet.Cells.Item(1, $i + 1).Value = $fileHeaders[$i]
}

# Open-ended (placeholder) list validations on "format" (D) and "data_type" (E)
$formatValidation = $fileSheet.Range("D2:D1048576").Validation
$formatValidation.Add(3, 1, 1, '""')
$formatValidation.ShowInput = $false
$formatValidation.ShowError = $false

$dataTypeValidation = $fileSheet.Range("E2:E1048576").Validation
$dataTypeValidation.Add(3, 1, 1, '""')
$dataTypeValidation.ShowInput = $false
$dataTypeValidation.ShowError = $false

# ---------------------------------------------------------------------------
# 2) "FileMetadata" worksheet - appended after "File"
# ---------------------------------------------------------------------------
$fileMetadataSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $fileSheet)
$fileMetadataSheet.Name = "FileMetadata"

$fileMetadataHeaders = @(
    "code",
    "display",
    "value_code",
    "value_display",
    "id",
    "external_id"
)

for ($i = 0; $i -lt $fileMetadataHeaders.Length; $i++) {
    $fileMetadataSheet.Cells.Item(1, $i + 1).Value = $fileMetadataHeaders[$i]
}

# Leave the first sheet ("Thing") selected/active, same as the original file.
$wb.Worksheets.Item(1).Activate()
